$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previously-blank, specially-formatted placeholder row (old row 2) is
# removed outright; every row below it shifts up by one (old row 3 becomes
# the new row 2, old row 4 becomes the new row 3, ..., old row 44 becomes
# the new row 43 - this is why the sheet's last row goes from 44 to 43).
$ws.Rows(2).Delete()

# Populate the (now) row 2 with the new group / absolute-position pair.
# Column B is written first so the shared-string table ends up in the same
# order as the target file (index 2 = "NZ_CP021201.1-1190942",
# index 3 = "te-09-0932-02").
$ws.Range("B2").Value = "NZ_CP021201.1-1190942"
$ws.Range("A2").Value = "te-09-0932-02"

# The new grouping label in column A is shown in the same bold header font
# used by the "Grouping"/"Absolute position" titles in row 1.
$ws.Range("A2").Font.Bold = $true
